$d = $word.ActiveDocument
$x = $d.ThisDoesNotExist12345
Write-Host "x=" $x
